# Weekly update for "Fruta / hortaliza" - Alcachofa, Mapocho Venta Directa de Santiago
# Applies the updated values for rows 2-15 (Fecha, Volumen, Precio minimo/maximo/promedio,
# Unidad de comercializacion, Origen, Precio $/Kg, Kg o Unidades).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44453
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 50
$ws.Range("K2").Value = 12000
$ws.Range("L2").Value = 12000
$ws.Range("M2").Value = 12000
$ws.Range("N2").Value = "$/caja 30 unidades"
$ws.Range("O2").Value = "Provincia de Limarí"
$ws.Range("P2").Value = 400
$ws.Range("Q2").Value = 30
$ws.Range("D3").Value = 44425
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 35
$ws.Range("K3").Value = 14000
$ws.Range("L3").Value = 14000
$ws.Range("M3").Value = 14000
$ws.Range("N3").Value = "$/caja 30 unidades"
$ws.Range("O3").Value = "Provincia de Limarí"
$ws.Range("P3").Value = 467
$ws.Range("Q3").Value = 30
$ws.Range("D4").Value = 44474
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 45
$ws.Range("K4").Value = 10000
$ws.Range("L4").Value = 10000
$ws.Range("M4").Value = 10000
$ws.Range("N4").Value = "$/caja 30 unidades"
$ws.Range("O4").Value = "Provincia de Limarí"
$ws.Range("P4").Value = 333
$ws.Range("Q4").Value = 30
$ws.Range("D5").Value = 44418
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 30
$ws.Range("K5").Value = 15000
$ws.Range("L5").Value = 15000
$ws.Range("M5").Value = 15000
$ws.Range("N5").Value = "$/caja 30 unidades"
$ws.Range("O5").Value = "Provincia de Limarí"
$ws.Range("P5").Value = 500
$ws.Range("Q5").Value = 30
$ws.Range("D6").Value = 44432
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 25
$ws.Range("K6").Value = 14000
$ws.Range("L6").Value = 14000
$ws.Range("M6").Value = 14000
$ws.Range("N6").Value = "$/caja 30 unidades"
$ws.Range("O6").Value = "Provincia del Elquí"
$ws.Range("P6").Value = 467
$ws.Range("Q6").Value = 30
$ws.Range("D7").Value = 44435
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 25
$ws.Range("K7").Value = 14000
$ws.Range("L7").Value = 14000
$ws.Range("M7").Value = 14000
$ws.Range("N7").Value = "$/caja 30 unidades"
$ws.Range("O7").Value = "Provincia de Limarí"
$ws.Range("P7").Value = 467
$ws.Range("Q7").Value = 30
$ws.Range("D8").Value = 44435
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 25
$ws.Range("K8").Value = 14000
$ws.Range("L8").Value = 14000
$ws.Range("M8").Value = 14000
$ws.Range("N8").Value = "$/caja 30 unidades"
$ws.Range("O8").Value = "Provincia del Elquí"
$ws.Range("P8").Value = 467
$ws.Range("Q8").Value = 30
$ws.Range("D9").Value = 44460
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 45
$ws.Range("K9").Value = 13000
$ws.Range("L9").Value = 13000
$ws.Range("M9").Value = 13000
$ws.Range("N9").Value = "$/caja 30 unidades"
$ws.Range("O9").Value = "Provincia de Limarí"
$ws.Range("P9").Value = 433
$ws.Range("Q9").Value = 30
$ws.Range("D10").Value = 44449
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 45
$ws.Range("K10").Value = 12000
$ws.Range("L10").Value = 12000
$ws.Range("M10").Value = 12000
$ws.Range("N10").Value = "$/caja 30 unidades"
$ws.Range("O10").Value = "Provincia de Limarí"
$ws.Range("P10").Value = 400
$ws.Range("Q10").Value = 30
$ws.Range("D11").Value = 44841
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 45
$ws.Range("K11").Value = 12000
$ws.Range("L11").Value = 12000
$ws.Range("M11").Value = 12000
$ws.Range("N11").Value = "$/caja 30 unidades"
$ws.Range("O11").Value = "Provincia de Limarí"
$ws.Range("P11").Value = 400
$ws.Range("Q11").Value = 30
$ws.Range("D12").Value = 44841
$ws.Range("I12").Value = "Segunda"
$ws.Range("J12").Value = 45
$ws.Range("K12").Value = 10000
$ws.Range("L12").Value = 10000
$ws.Range("M12").Value = 10000
$ws.Range("N12").Value = "$/caja 40 unidades"
$ws.Range("O12").Value = "Provincia de Limarí"
$ws.Range("P12").Value = 250
$ws.Range("Q12").Value = 40
$ws.Range("D13").Value = 44421
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 25
$ws.Range("K13").Value = 15000
$ws.Range("L13").Value = 16000
$ws.Range("M13").Value = 15400
$ws.Range("N13").Value = "$/caja 30 unidades"
$ws.Range("O13").Value = "Provincia del Elquí"
$ws.Range("P13").Value = 513
$ws.Range("Q13").Value = 30
$ws.Range("D14").Value = 44376
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 25
$ws.Range("K14").Value = 18000
$ws.Range("L14").Value = 18000
$ws.Range("M14").Value = 18000
$ws.Range("N14").Value = "$/caja 30 unidades"
$ws.Range("O14").Value = "Provincia de Limarí"
$ws.Range("P14").Value = 600
$ws.Range("Q14").Value = 30
$ws.Range("D15").Value = 44446
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 25
$ws.Range("K15").Value = 14000
$ws.Range("L15").Value = 14000
$ws.Range("M15").Value = 14000
$ws.Range("N15").Value = "$/caja 30 unidades"
$ws.Range("O15").Value = "Provincia de Limarí"
$ws.Range("P15").Value = 467
$ws.Range("Q15").Value = 30
